$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used data row
$lastRow = $ws.Range("A1").End(-4121).Row  # xlDown = -4121

# Add new headers in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style/format (bold, border, centered) from A1 to AD1:AF1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in team record values for each data row (2 through lastRow)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 84   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 78   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
